$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case fixes for "de"/"el"/"la"/"los" connector words ---
$ws.Range("B4").Value = "Mazapa De Madero"
$ws.Range("A10").Value = "Ciudad De México"
$ws.Range("A25").Value = "Estado De México"
$ws.Range("B25").Value = "Ecatepec De Morelos"
$ws.Range("B30").Value = "Naucalpan De Juárez"
$ws.Range("B34").Value = "Apaseo El Alto"
$ws.Range("B36").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B40").Value = "Acapulco De Juárez"
$ws.Range("B42").Value = "Ayutla De Los Libres"
$ws.Range("B45").Value = "Huitzuco De Los Figueroa"
$ws.Range("B50").Value = "Técpan De Galeana"
$ws.Range("B51").Value = "Tlapa De Comonfort"
$ws.Range("B56").Value = "Huasca De Ocampo"
$ws.Range("B59").Value = "Mixquiahuala De Juárez"
$ws.Range("B60").Value = "Omitlán De Juárez"
$ws.Range("B61").Value = "Pachuca De Soto"
$ws.Range("B66").Value = "Tulancingo De Bravo"
$ws.Range("B69").Value = "Autlán De Navarro"
$ws.Range("B72").Value = "Teocuitatlán De Corona"
$ws.Range("B91").Value = "Zimatlán De Álvarez"
$ws.Range("B99").Value = "Huehuetlán El Chico"
$ws.Range("B101").Value = "Izúcar De Matamoros"
$ws.Range("B104").Value = "San Nicolás De Los Ranchos"
$ws.Range("B112").Value = "Jalpan De Serra"
$ws.Range("B113").Value = "Landa De Matamoros"
$ws.Range("B128").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B138").Value = "Hueyapan De Ocampo"
$ws.Range("B139").Value = "Ignacio De La Llave"

# --- Tiny floating point recalculation tweaks ---
$ws.Range("D53").Value = 0.09740259740259739
$ws.Range("D111").Value = 0.09090909090909093

# --- Remove trailing source/footer rows (154-158) ---
$ws.Range("A154:A158").EntireRow.Delete() | Out-Null
